$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows for species records at rows 6, 7 and 8 got re-sorted/rotated.
# New row 6 = old row 7, new row 7 = old row 8, new row 8 = old row 6.
# Capture the current ("before") values of the columns that actually change
# for each of the three rows, then write them back in rotated order.

$cols = @("A","B","D","E","F","G","H","P","Q","R","AW","AX")

$old6 = @{}
$old7 = @{}
$old8 = @{}
foreach ($col in $cols) {
    $old6[$col] = $ws.Range($col + "6").Value2
    $old7[$col] = $ws.Range($col + "7").Value2
    $old8[$col] = $ws.Range($col + "8").Value2
}

foreach ($col in $cols) {
    $ws.Range($col + "6").Value2 = $old7[$col]
    $ws.Range($col + "7").Value2 = $old8[$col]
    $ws.Range($col + "8").Value2 = $old6[$col]
}
